$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4: Literature Society IITJ Website
$ws.Range("E4").Value = 0.85
$ws.Range("F4").Value = 62.32

# Row 5: CloudPhysician's Vital Extraction Challenge
$ws.Range("E5").Value = 0.8
$ws.Range("F5").Value = 58.66

# Row 6: FaceNet Implementation
$ws.Range("E6").Value = 0.8
$ws.Range("F6").Value = 58.66

# Row 10: Cloudphysician's Vital Extraction Challenge
$ws.Range("E10").Value = 0.8
$ws.Range("F10").Value = 64

# Row 11: Website for the Literature Society of the college
$ws.Range("E11").Value = 0.85
$ws.Range("F11").Value = 65.14
